$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.207.94'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +14.87%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.676.61'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +9.17%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.94%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +9.08%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9978'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +3.70%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3717'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.67%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3435'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +7.84%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '48.25'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +18.77%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.182'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +8.17%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07277'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +7.06%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9994'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.57'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +10.12%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.077'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +7.29%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.748'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +6.35%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.677.73'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +9.74%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001109'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +6.17%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9975'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.69%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +10.35%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '81.71'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +13.07%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.45'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +9.85%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.120'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +7.58%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +5.80%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.206.75'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +14.25%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.398'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.362'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -8.84%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +21.07%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '151.56'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.55'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +10.50%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.861.12'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +9.50%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '127.03'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +7.44%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.442'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +24.46%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.035'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9900'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +16.38%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.738'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +15.12%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08437'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +5.58%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +16.70%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Hedera'
$ws.Range('B38').Style = "Normal"
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C38').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06405'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +9.34%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('B39').Style = "Normal"
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C39').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.358'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +8.04%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.883'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +15.21%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.291'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +7.16%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.02340'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +11.03%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.2106'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +9.82%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6144'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +13.18%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9964'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '13.24'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +6.23%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('B47').Style = "Normal"
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C47').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.801'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5953'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +9.24%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '127.51'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +5.08%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.016'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +8.16%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07131'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +8.31%  '
$ws.Range('E51').Style = "Normal"
